# Generate Report for Handback
#
# This script mutates the "localization-status" workbook so that, for both
# the zh-cn and de-de handoff rows, the "Latest Target File" / "Latest
# Handback File" / "Latest Handback DateTime" columns are populated (the
# handback has now completed and is in sync with en-US), and the Status
# columns flip from "Ready for handoff" to "Handed back: in sync with
# en-US".

$wb = $excel.ActiveWorkbook

# ---------------------------------------------------------------------
# 1. Status text: "Ready for handoff" -> "Handed back: in sync with en-US"
#    Every cell that used to read "Ready for handoff" gets the new text
#    (Overview sheet status columns + per-locale Status column).
# ---------------------------------------------------------------------
$newStatus = "Handed back: in sync with en-US"

$wsOverview = $wb.Worksheets.Item("Overview")
$wsOverview.Range("E2").Value = $newStatus
$wsOverview.Range("F2").Value = $newStatus
$wsOverview.Range("E3").Value = $newStatus
$wsOverview.Range("F3").Value = $newStatus

$wsZh = $wb.Worksheets.Item("zh-cn")
$wsZh.Range("C2").Value = $newStatus
$wsZh.Range("C3").Value = $newStatus

$wsDe = $wb.Worksheets.Item("de-de")
$wsDe.Range("C2").Value = $newStatus
$wsDe.Range("C3").Value = $newStatus

# ---------------------------------------------------------------------
# 2. Per-row handback info: zh-cn sheet.
#    I = Latest Target File (hyperlink to the source .md, like column A)
#    J = Latest Handback File (the generated handback xliff file name)
#    K = Latest Handback DateTime (already populated -> refresh value)
# ---------------------------------------------------------------------
$mdUrl1 = "https://github.com/OpenLocalizationTestOrg/ol-test0/blob/d96997bc32300e25d4647a9f1eb9e22ef4ee286e/e2e/43a2342d-3ef4-4fdd-898d-e284617deb68.md"
$mdUrl2 = "https://github.com/OpenLocalizationTestOrg/ol-test0/blob/d96997bc32300e25d4647a9f1eb9e22ef4ee286e/e2e/bdd835bf-7442-4ff3-8b76-10514be9fdd3.md"
$mdName1 = "43a2342d-3ef4-4fdd-898d-e284617deb68.md"
$mdName2 = "bdd835bf-7442-4ff3-8b76-10514be9fdd3.md"

# Re-create the sheet hyperlink collection in the order:
#   A2 (existing), I2 (new), A3 (existing), I3 (new)
# so the relationship ids line up the way Excel would assign them when the
# two new hyperlinks are inserted after the original A2 link and before the
# original A3 link is re-pointed.
$wsZh.Hyperlinks.Delete()

$wsZh.Range("J2").Value = "43a2342d-3ef4-4fdd-898d-e284617deb68.3713d90633c05ce8a8e6b14d45f2e06a436cebad.zh-cn.xlf"
$wsZh.Range("K2").Value = "2016-08-19 02:44:35"

$wsZh.Range("J3").Value = "bdd835bf-7442-4ff3-8b76-10514be9fdd3.564fb869d5dd82b2b74ef4e32f68e16104ef02b4.zh-cn.xlf"
$wsZh.Range("K3").Value = "2016-08-19 02:44:35"

$wsZh.Hyperlinks.Add($wsZh.Range("A2"), $mdUrl1, "", "", $mdName1)
$wsZh.Hyperlinks.Add($wsZh.Range("I2"), $mdUrl1, "", "", $mdName1)
$wsZh.Hyperlinks.Add($wsZh.Range("A3"), $mdUrl2, "", "", $mdName2)
$wsZh.Hyperlinks.Add($wsZh.Range("I3"), $mdUrl2, "", "", $mdName2)

$wsZh.Range("I2").Font.Underline = 2
$wsZh.Range("I2").Font.Color = 15570276
$wsZh.Range("I3").Font.Underline = 2
$wsZh.Range("I3").Font.Color = 15570276

# ---------------------------------------------------------------------
# 3. Per-row handback info: de-de sheet (same shape as zh-cn, but K gets a
#    distinct handback timestamp).
# ---------------------------------------------------------------------
$wsDe.Hyperlinks.Delete()

$wsDe.Range("J2").Value = "43a2342d-3ef4-4fdd-898d-e284617deb68.3713d90633c05ce8a8e6b14d45f2e06a436cebad.de-de.xlf"
$wsDe.Range("K2").Value = "2016-08-19 02:44:43"

$wsDe.Range("J3").Value = "bdd835bf-7442-4ff3-8b76-10514be9fdd3.564fb869d5dd82b2b74ef4e32f68e16104ef02b4.de-de.xlf"
$wsDe.Range("K3").Value = "2016-08-19 02:44:43"

$wsDe.Hyperlinks.Add($wsDe.Range("A2"), $mdUrl1, "", "", $mdName1)
$wsDe.Hyperlinks.Add($wsDe.Range("I2"), $mdUrl1, "", "", $mdName1)
$wsDe.Hyperlinks.Add($wsDe.Range("A3"), $mdUrl2, "", "", $mdName2)
$wsDe.Hyperlinks.Add($wsDe.Range("I3"), $mdUrl2, "", "", $mdName2)

$wsDe.Range("I2").Font.Underline = 2
$wsDe.Range("I2").Font.Color = 15570276
$wsDe.Range("I3").Font.Underline = 2
$wsDe.Range("I3").Font.Color = 15570276

# ---------------------------------------------------------------------
# 4. Column widths: the Status columns (Overview E/F, per-locale C) and the
#    newly-populated Latest Target File / Latest Handback File columns
#    (per-locale I/J) widen to fit the longer text now shown in them.
# ---------------------------------------------------------------------
$wideStatus = 29.09   # -> stored width ~29.98 (character-width rounding)
$wideFile = 39.17     # -> stored width 40 (matches other 40-wide columns)

$wsOverview.Columns.Item(5).ColumnWidth = $wideStatus
$wsOverview.Columns.Item(6).ColumnWidth = $wideStatus

$wsZh.Columns.Item(3).ColumnWidth = $wideStatus
$wsZh.Columns.Item(9).ColumnWidth = $wideFile
$wsZh.Columns.Item(10).ColumnWidth = $wideFile

$wsDe.Columns.Item(3).ColumnWidth = $wideStatus
$wsDe.Columns.Item(9).ColumnWidth = $wideFile
$wsDe.Columns.Item(10).ColumnWidth = $wideFile
